$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = "Boots"
$ws.Range("A15").Value = "King"
$ws.Range("B14").Value = "The-chain.nl/story/characters/boots"
$ws.Range("B15").Value = "The-chain.nl/story/characters/king"

$ws.Range("B29").Select()
